# Update cryptos list values per upstream diff (row-level value refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns that we are about to rewrite to stay text-typed
# (their source values include things like "26.441.73" and "0.000008703" which
# must round-trip as literal text, not be coerced into numbers).
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '26.441.73'
$ws.Range("E2").Value = '  -2.81%  '

# Row 3
$ws.Range("D3").Value = '1.774.66'
$ws.Range("E3").Value = '  -1.67%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("B5").Value = 'USDC'
$ws.Range("C5").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D5").Value = '1.004'
$ws.Range("E5").Value = '  +0.04%  '

# Row 6
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '305.67'
$ws.Range("E6").Value = '  -1.66%  '

# Row 7
$ws.Range("D7").Value = '0.4281'
$ws.Range("E7").Value = '  +1.65%  '

# Row 8
$ws.Range("D8").Value = '0.3628'
$ws.Range("E8").Value = '  +2.08%  '

# Row 9
$ws.Range("D9").Value = '0.07154'
$ws.Range("E9").Value = '  +0.41%  '

# Row 10
$ws.Range("D10").Value = '0.8416'
$ws.Range("E10").Value = '  -0.58%  '

# Row 11
$ws.Range("D11").Value = '20.48'
$ws.Range("E11").Value = '  +1.47%  '

# Row 12
$ws.Range("D12").Value = '1.789.32'
$ws.Range("E12").Value = '  -6.75%  '

# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = '6.449'
$ws.Range("E13").Value = '  +1.39%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.249'
$ws.Range("E14").Value = '  -1.51%  '

# Row 15
$ws.Range("D15").Value = '0.06887'
$ws.Range("E15").Value = '  -0.48%  '

# Row 16
$ws.Range("D16").Value = '1.008'
$ws.Range("E16").Value = '  +0.18%  '

# Row 17
$ws.Range("D17").Value = '78.74'
$ws.Range("E17").Value = '  -2.40%  '

# Row 18
$ws.Range("D18").Value = '0.000008703'
$ws.Range("E18").Value = '  -1.01%  '

# Row 19
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  -0.38%  '

# Row 20
$ws.Range("D20").Value = '14.91'
$ws.Range("E20").Value = '  -1.11%  '

# Row 21
$ws.Range("D21").Value = '26.456.10'
$ws.Range("E21").Value = '  -4.45%  '

# Row 22
$ws.Range("D22").Value = '5.103'
$ws.Range("E22").Value = '  +0.57%  '

# Row 23
$ws.Range("D23").Value = '11.10'
$ws.Range("E23").Value = '  +2.27%  '

# Row 24
$ws.Range("D24").Value = '2.023.51'
$ws.Range("E24").Value = '  -4.29%  '

# Row 25
$ws.Range("D25").Value = '152.42'
$ws.Range("E25").Value = '  -0.63%  '

# Row 26
$ws.Range("D26").Value = '1.868'
$ws.Range("E26").Value = '  -4.64%  '

# Row 27
$ws.Range("D27").Value = '18.02'
$ws.Range("E27").Value = '  -1.21%  '

# Row 28
$ws.Range("D28").Value = '5.066'
$ws.Range("E28").Value = '  +0.23%  '

# Row 29
$ws.Range("D29").Value = '113.82'
$ws.Range("E29").Value = '  +0.75%  '

# Row 30
$ws.Range("D30").Value = '1.794'
$ws.Range("E30").Value = '  +4.69%  '

# Row 31
$ws.Range("D31").Value = '0.08889'
$ws.Range("E31").Value = '  -0.08%  '

# Row 32
$ws.Range("D32").Value = '0.7252'
$ws.Range("E32").Value = '  -2.27%  '

# Row 33
$ws.Range("D33").Value = '1.121'
$ws.Range("E33").Value = '  +1.55%  '

# Row 34
$ws.Range("D34").Value = '4.325'
$ws.Range("E34").Value = '  -3.09%  '

# Row 35
$ws.Range("B35").Value = 'Frax'
$ws.Range("C35").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D35").Value = '1.004'
$ws.Range("E35").Value = '  -0.06%  '

# Row 36
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '2.739'
$ws.Range("E36").Value = '  -6.93%  '

# Row 37
$ws.Range("D37").Value = '1.101'
$ws.Range("E37").Value = '  +2.81%  '

# Row 38
$ws.Range("D38").Value = '0.05138'
$ws.Range("E38").Value = '  -1.53%  '

# Row 39
$ws.Range("D39").Value = '0.01888'
$ws.Range("E39").Value = '  -0.55%  '

# Row 40
$ws.Range("D40").Value = '0.1612'
$ws.Range("E40").Value = '  -1.47%  '

# Row 41
$ws.Range("D41").Value = '0.4920'
$ws.Range("E41").Value = '  -1.09%  '

# Row 42
$ws.Range("D42").Value = '2.597'
$ws.Range("E42").Value = '  -4.82%  '

# Row 43
$ws.Range("D43").Value = '6.324'
$ws.Range("E43").Value = '  +0.49%  '

# Row 44
$ws.Range("D44").Value = '7.983'
$ws.Range("E44").Value = '  -2.52%  '

# Row 45
$ws.Range("D45").Value = '104.86'
$ws.Range("E45").Value = '  -0.09%  '

# Row 46
$ws.Range("D46").Value = '1.004'
$ws.Range("E46").Value = '  +0.08%  '

# Row 47
$ws.Range("D47").Value = '10.14'
$ws.Range("E47").Value = '  -0.97%  '

# Row 48
$ws.Range("D48").Value = '1.632'
$ws.Range("E48").Value = '  +2.34%  '

# Row 49
$ws.Range("D49").Value = '0.06194'
$ws.Range("E49").Value = '  -3.01%  '

# Row 50
$ws.Range("D50").Value = '0.4461'
$ws.Range("E50").Value = '  -2.27%  '

# Row 51
$ws.Range("D51").Value = '1.705'
$ws.Range("E51").Value = '  +1.96%  '
